$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for "Acelga" at Feria Lagunitas de
# Puerto Montt. It is inserted as a new row 146 (most recent first), which
# pushes every existing record from row 146 down one row (146->147, ...,
# 169->170) and extends the used range from A1:R169 to A1:R170.
$ws.Rows.Item(146).Insert()

$ws.Cells.Item(146, 1).Value = 4
$ws.Cells.Item(146, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(146, 3).Value = "Los Lagos"
$ws.Cells.Item(146, 4).Value = 44694
$ws.Cells.Item(146, 5).Value = 10
$ws.Cells.Item(146, 6).Value = 100112009
$ws.Cells.Item(146, 7).Value = "Acelga"
$ws.Cells.Item(146, 8).Value = "Sin especificar"
$ws.Cells.Item(146, 9).Value = "Primera"
$ws.Cells.Item(146, 10).Value = 90
$ws.Cells.Item(146, 11).Value = 10000
$ws.Cells.Item(146, 12).Value = 10000
$ws.Cells.Item(146, 13).Value = 10000
$ws.Cells.Item(146, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(146, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(146, 16).Value = 833
$ws.Cells.Item(146, 17).Value = 12
$ws.Cells.Item(146, 18).Value = "Hortaliza"
